$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.862
$ws.Range("B4").Value = 7.043000000000001

$ws.Range("B5").Value = 6.38

$ws.Range("A6").Value = -21.108

$ws.Range("A7").Value = -21.038

$ws.Range("B8").Value = 6.161

$ws.Range("A16").Value = -21.071
$ws.Range("B16").Value = 6.528

$ws.Range("A20").Value = -21.98

$ws.Range("B22").Value = 6.625999999999999
